$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2024 commemorative coin row: mintage figure became known (1.000.000)
$ws.Range("F22").Value = "1.000.000"
$ws.Range("G22").Value = 1

# Leave the cursor where the author last left it when saving
$ws.Range("F26").Select()
